# Minimal-case test fixture touch-up.
#
# The header in column C was renamed from "replicate_number" to
# "replicatenumber", the ad-hoc formatting that had been stamped onto the
# data rows (B2:H6) is cleared back to the workbook's Normal style, and the
# active selection is left on C2 (matching the resaved file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "replicate_number" header to "replicatenumber".
$ws.Range("C1").Value = "replicatenumber"

# Drop the stray formatting that had been applied to the data cells in
# columns B:H (rows 2-6), restoring them to the default/unstyled look.
$ws.Range("B2:H6").Style = "Normal"

# Leave the selection where the author left it when the file was saved.
$ws.Range("C2").Select() | Out-Null
